$d = $word.ActiveDocument

# Locate the paragraph that contains the sentence being edited by anchoring on a
# unique, unaffected piece of text ("задержание вооруженных преступников").
$anchor = $d.Content
$anchor.Find.Execute("задержание вооруженных преступников", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $anchor.Paragraphs(1)
$rng = $targetPara.Range

# Replace the whole paragraph's OOXML in one shot so the run layout (and the
# now-stale gramStart/gramEnd proofErr markers around the deleted sentence)
# comes out exactly as a human edit in Word would leave it: the sentence
# "Входил в состав штурмовой группы" is gone, the trailing period from that
# sentence collapses onto the previous run, and the following sentence
# (" В дополнении к основным обязанностям...") becomes its own run.
$newParaXml = '<w:p w:rsidR="00D630CA" w:rsidRPr="00B35981" w:rsidRDefault="00D630CA" w:rsidP="00D630CA"><w:pPr><w:ind w:firstLine="709"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00596A33"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Добрый день члены комиссии, меня зовут Дарин Сергей, мне 36 лет, я из города Тольятти. На данный момент я военный пенсионер. На предыдущем месте работы проработал</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> с 2009</w:t></w:r><w:r w:rsidRPr="00596A33"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> по конец 2022</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> года</w:t></w:r><w:r w:rsidRPr="00596A33"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidRPr="00B35981"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Род деятельности – служебные командировки в район СКР, участие в проведении КТО, задержание вооруженных преступников</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> В дополнении к основным обязанностям, я был внештатным санинструктором и отвечал за подготовку личного состава в области тактической медицины (проводил теоретические и практические занятия).</w:t></w:r></w:p>'

$rng.InsertXML($newParaXml)
